$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. "Stop Date: 9/30/2015"  ->  "Stop Date: "  (drop the date, keep the
#    trailing space so the run ends up with xml:space="preserve").
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("Stop Date: 9/30/2015", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Stop Date: ", 2) | Out-Null

# ---------------------------------------------------------------------------
# 2. After the "1. Planning Phase" line, insert four new detail lines
#    (each starting with a tab) that capture the planned/actual
#    start & stop dates for this phase, ending with the relocated
#    "_GoBack" bookmark.
# ---------------------------------------------------------------------------
$xmlHeader = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' + `
  '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
  '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
  '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>'
$xmlFooter = '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$planningPhasePara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "1. Planning Phase*") {
        $planningPhasePara = $p
        break
    }
}

$planningPhasePara.Range.InsertParagraphAfter() | Out-Null

$newLinesXml = $xmlHeader + `
  '<w:p><w:r><w:tab/><w:t>Planned Start: 9/1/2015</w:t></w:r></w:p>' + `
  '<w:p><w:r><w:tab/><w:t>Start: 9/1/2015</w:t></w:r></w:p>' + `
  '<w:p><w:r><w:tab/><w:t>Planned Stop: 9/30/2015</w:t></w:r></w:p>' + `
  '<w:p><w:r><w:tab/><w:t xml:space="preserve">Stop: </w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>' + `
  $xmlFooter

$placeholderPara = $planningPhasePara.Next()
$placeholderPara.Range.InsertXML($newLinesXml) | Out-Null

# ---------------------------------------------------------------------------
# 3. Collapse the old "9. " / proofErr / "…" / bookmark paragraph down to a
#    single plain run "9. …" (the bookmark now lives on the new "Stop: "
#    line above, so it is not repeated here).
# ---------------------------------------------------------------------------
$nineDotsPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "9.*") {
        $nineDotsPara = $p
        break
    }
}

$nineDotsXml = $xmlHeader + '<w:p><w:r><w:t>9. &#8230;</w:t></w:r></w:p>' + $xmlFooter
$nineDotsPara.Range.InsertXML($nineDotsXml) | Out-Null
